$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 13; $row++) {
    # Column H: PERIOD TO EXPIRE -> decrement by 1
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value2 - 1

    # Column I: LAST UPDATE -> set to 04-Nov-2025, keep as plain text
    # (leading apostrophe forces Excel to store it as text instead of
    # auto-converting the date-looking string into a date serial value)
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Formula = "'04-Nov-2025"
}
